$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated measurement data for columns B, E, H (rows 2-6) with new plunger data
$ws.Range("B2").Value = 298.8
$ws.Range("E2").Value = 296.8
$ws.Range("H2").Value = 296.7

$ws.Range("B3").Value = 299.7
$ws.Range("E3").Value = 297.3
$ws.Range("H3").Value = 297

$ws.Range("B4").Value = 300.5
$ws.Range("E4").Value = 297.9
$ws.Range("H4").Value = 296.1

$ws.Range("B5").Value = 300.1
$ws.Range("E5").Value = 297.8
$ws.Range("H5").Value = 296.3

$ws.Range("B6").Value = 300.2
$ws.Range("E6").Value = 297.8
$ws.Range("H6").Value = 297
